# Apply the "Updated cryptos list" refresh: new prices and 1h volume
# percentages scraped for each coin, plus a ranking swap for three pairs
# of rows (NEARProtocol/EthereumClassic and Fetch.AI/Kaspa/Stacks).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds values that look numeric (e.g. "0.632")
# but must stay plain text, matching the sheet's existing inline-string
# cells (no numeric coercion, no thousands separators as numbers).
# Forcing NumberFormat to Text before the write - then clearing the
# format again so no stray style survives - keeps the cell a string
# without leaving any formatting residue behind.
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue "D2" "69.858.90"
$ws.Range("E2").Value = "  -1.28%  "
Set-TextValue "D3" "3.577.95"
$ws.Range("E3").Value = "  -1.29%  "
$ws.Range("E4").Value = "  -0.11%  "
Set-TextValue "D5" "577.40"
$ws.Range("E5").Value = "  -2.30%  "
Set-TextValue "D6" "189.77"
$ws.Range("E6").Value = "  -1.10%  "
Set-TextValue "D7" "0.632"
$ws.Range("E7").Value = "  -2.91%  "
Set-TextValue "D8" "3.572.91"
$ws.Range("E8").Value = "  -1.22%  "
$ws.Range("E9").Value = "  +0.01%  "
Set-TextValue "D10" "0.179"
$ws.Range("E10").Value = "  -0.45%  "
Set-TextValue "D11" "0.661"
$ws.Range("E11").Value = "  -0.47%  "
Set-TextValue "D12" "55.77"
$ws.Range("E12").Value = "  -4.50%  "
Set-TextValue "D13" "0.0000305"
$ws.Range("E13").Value = "  +3.10%  "
Set-TextValue "D14" "9.65"
$ws.Range("E14").Value = "  -1.55%  "
Set-TextValue "D15" "4.151.23"
$ws.Range("E15").Value = "  -1.26%  "
Set-TextValue "D16" "19.86"
$ws.Range("E16").Value = "  +2.25%  "
Set-TextValue "D17" "3.570.84"
$ws.Range("E17").Value = "  -1.59%  "
Set-TextValue "D18" "69.796.59"
$ws.Range("E18").Value = "  -1.40%  "
Set-TextValue "D19" "12.65"
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("E20").Value = "  +0.00%  "
$ws.Range("E21").Value = "  -1.27%  "
Set-TextValue "D22" "476.72"
$ws.Range("E22").Value = "  -4.41%  "
Set-TextValue "D23" "19.40"
$ws.Range("E23").Value = "  +12.27%  "
Set-TextValue "D24" "5.05"
$ws.Range("E24").Value = "  -6.19%  "
Set-TextValue "D25" "95.78"
$ws.Range("E25").Value = "  +5.07%  "
$ws.Range("E26").Value = "  -2.92%  "
Set-TextValue "D27" "3.00"
$ws.Range("E27").Value = "  -3.98%  "
Set-TextValue "D28" "11.04"
$ws.Range("E28").Value = "  -1.52%  "
Set-TextValue "D29" "9.38"
$ws.Range("E29").Value = "  -0.87%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue "D30" "32.35"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D31" "7.70"
$ws.Range("E31").Value = "  +1.47%  "
Set-TextValue "D32" "12.23"
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("E33").Value = "  +0.64%  "
Set-TextValue "D34" "66.33"
$ws.Range("E34").Value = "  +1.77%  "
Set-TextValue "D35" "582.42"
$ws.Range("E35").Value = "  -6.88%  "
Set-TextValue "D36" "38.93"
$ws.Range("E36").Value = "  +1.90%  "
Set-TextValue "D38" "0.0₃0803"
$ws.Range("E38").Value = "  -3.19%  "
Set-TextValue "D39" "0.396"
$ws.Range("E39").Value = "  -3.17%  "
Set-TextValue "D40" "3.26"
$ws.Range("E40").Value = "  +19.70%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D41" "0.138"
$ws.Range("E41").Value = "  -6.43%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
Set-TextValue "D42" "3.46"
$ws.Range("E42").Value = "  -5.48%  "
$ws.Range("B43").Value = "Fetch.AI"
$ws.Range("C43").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue "D43" "2.87"
$ws.Range("E43").Value = "  +6.83%  "
Set-TextValue "D44" "3.238.12"
$ws.Range("E44").Value = "  -2.61%  "
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("E46").Value = "  -1.12%  "
Set-TextValue "D47" "3.36"
$ws.Range("E47").Value = "  +0.15%  "
Set-TextValue "D48" "9.31"
$ws.Range("E48").Value = "  +1.51%  "
$ws.Range("E49").Value = "  +0.08%  "
Set-TextValue "D50" "0.998"
$ws.Range("E50").Value = "  -0.24%  "
Set-TextValue "D51" "3.13"
$ws.Range("E51").Value = "  -5.34%  "
